$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BasicUser")

# Append the new "Anaïs Rouvière" user as row 11
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Anaïs"
$ws.Range("C11").Value = "Rouvière"
$ws.Range("D11").Value = "arouviere"
$ws.Range("E11").Value = "anais.rouviere@kobalt.fr"
$ws.Range("F11").Value = "kobalt"
$ws.Range("G11").Value = 2
$ws.Range("H11").Value = 2
$ws.Range("I11").Value = "fr"

$ws.Hyperlinks.Add($ws.Range("E11"), "mailto:anais.rouviere@kobalt.fr", "", "", "anais.rouviere@kobalt.fr")
$ws.Range("E11").Font.Underline = $false
$ws.Range("E11").Font.Color = 16711680
$ws.Range("E11").Font.Name = "Arial"
$ws.Range("E11").Font.Size = 10

$ws.Range("B2").Select() | Out-Null
